$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price and volume figures (symbol list refresh).
# D and E columns are stored as text in the workbook (e.g. "321.07", "-3.40%"),
# so force Text number format before assigning to avoid Excel auto-converting
# these numeric-looking strings into actual numbers/percentages.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D2").Value = "321.07"
$ws.Range("E2").Value = "-3.40%"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D3").Value = "42.88"
$ws.Range("E3").Value = "-5.49%"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D4").Value = "5.209"
$ws.Range("E4").Value = "-5.64%"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08180"
$ws.Range("E5").Value = "-3.46%"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D6").Value = "4.329"
$ws.Range("E6").Value = "-2.56%"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D7").Value = "1.804"
$ws.Range("E7").Value = "-12.81%"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9488"
$ws.Range("E8").Value = "-3.86%"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1123"
$ws.Range("E9").Value = "-4.46%"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1861"
$ws.Range("E10").Value = "-3.23%"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09425"
$ws.Range("E11").Value = "-4.19%"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04611"
$ws.Range("E12").Value = "-2.00%"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D13").Value = "7.453"
$ws.Range("E13").Value = "-21.56%"

# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.13%"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001295"
$ws.Range("E15").Value = "0.85%"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005855"
$ws.Range("E16").Value = "-0.78%"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D17").Value = "3.359"
$ws.Range("E17").Value = "-0.84%"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D18").Value = "2.544"
$ws.Range("E18").Value = "0.07%"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3362"
$ws.Range("E19").Value = "0.95%"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1389"
$ws.Range("E20").Value = "0.22%"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2551"
$ws.Range("E21").Value = "-0.01%"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04138"
$ws.Range("E22").Value = "-0.57%"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001249"
$ws.Range("E23").Value = "-4.21%"

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-6.97%"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0002981"
$ws.Range("E26").Value = "-0.28%"

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02646"
$ws.Range("E38").Value = "-2.30%"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05556"
$ws.Range("E39").Value = "-3.44%"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D40").Value = "0.008153"
$ws.Range("E40").Value = "4.48%"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1399"
$ws.Range("E41").Value = "-2.37%"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006563"
$ws.Range("E42").Value = "-12.06%"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002089"
$ws.Range("E43").Value = "-3.25%"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007585"
$ws.Range("E44").Value = "-6.03%"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3202"
$ws.Range("E45").Value = "-9.91%"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006850"
$ws.Range("E46").Value = "-3.14%"

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.28%"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003337"
$ws.Range("E48").Value = "-3.21%"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004103"
$ws.Range("E49").Value = "15.82%"

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.28%"

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.28%"
